$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Professional summary paragraph: plain-text substitution
#    "affecting all Black and Asian-American voters" -> "affecting 50M voters"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic demographic coding errors affecting 50M voters, developed",
    2)

# ---------------------------------------------------------------------
# 2) "Impact:" paragraph: plain-text substitution
#    "affecting all Black and Asian-American voters," -> "affecting 50M voters nationwide,"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved",
    2)

# ---------------------------------------------------------------------
# 3) Bullet point under "Partner - Siege Analytics": this run needs to be
#    split into three runs so that "50M" picks up bold + the accent color
#    used elsewhere in the document for emphasized figures, matching:
#      "...affecting " + [b, color=2C3E50]"50M" + " voters, developed..."
#    After steps 1 and 2 above, this is now the only remaining occurrence
#    of "all Black and Asian-American" in the document.
# ---------------------------------------------------------------------
$found = $d.Content
$found.Find.Execute("all Black and Asian-American", $true)

$target = $d.Range($found.Start, $found.End)
$target.Text = "50M"

$bolded = $d.Range($found.Start, $found.Start + 3)
$bolded.Font.Bold = 1
$bolded.Font.Color = 5258796
